# Auto-generated edit script: update cryptos list values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.983.29"
$ws.Range("E2").Value = "  -0.23%  "

$ws.Range("D3").Value = "1.825.87"
$ws.Range("E3").Value = "  -0.29%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9960"
$ws.Range("E4").Value = "  -0.32%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.72"
$ws.Range("E5").Value = "  +0.88%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6302"
$ws.Range("E6").Value = "  +0.52%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9985"
$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07500"
$ws.Range("E8").Value = "  -1.28%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2935"
$ws.Range("E9").Value = "  +0.53%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.02"
$ws.Range("E10").Value = "  +0.74%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07692"
$ws.Range("E11").Value = "  +0.67%  "

$ws.Range("D12").Value = "1.837.65"
$ws.Range("E12").Value = "  +0.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.989"
$ws.Range("E13").Value = "  +0.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6666"
$ws.Range("E14").Value = "  +0.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.96"
$ws.Range("E15").Value = "  +0.68%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009605"
$ws.Range("E16").Value = "  +2.15%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.042"
$ws.Range("E17").Value = "  +0.86%  "

$ws.Range("D18").Value = "29.003.11"
$ws.Range("E18").Value = "  +0.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.55"
$ws.Range("E19").Value = "  +1.86%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "225.47"
$ws.Range("E20").Value = "  +0.25%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9973"
$ws.Range("E21").Value = "  -0.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.135"
$ws.Range("E22").Value = "  -1.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9977"
$ws.Range("E23").Value = "  -0.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "160.23"
$ws.Range("E24").Value = "  -0.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1417"
$ws.Range("E25").Value = "  +3.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.494"
$ws.Range("E26").Value = "  +0.93%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.89"
$ws.Range("E27").Value = "  +0.25%  "

$ws.Range("E28").Value = "  +0.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.129"
$ws.Range("E29").Value = "  +1.74%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.054"
$ws.Range("E30").Value = "  +0.45%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05456"
$ws.Range("E31").Value = "  +4.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.198"
$ws.Range("E32").Value = "  +0.20%  "

$ws.Range("E33").Value = "  +0.16%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7430"
$ws.Range("E34").Value = "  +1.84%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.135"
$ws.Range("E35").Value = "  -1.46%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.640"
$ws.Range("E36").Value = "  +1.78%  "

$ws.Range("D37").Value = "1.238.92"
$ws.Range("E37").Value = "  -2.72%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.749"
$ws.Range("E38").Value = "  -0.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01778"
$ws.Range("E39").Value = "  -0.43%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.641"
$ws.Range("E40").Value = "  +2.04%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8994"
$ws.Range("E41").Value = "  +0.79%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9984"
$ws.Range("E42").Value = "  -0.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.24"
$ws.Range("E43").Value = "  -0.33%  "

$ws.Range("D44").Value = "1.976.08"
$ws.Range("E44").Value = "  +0.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000124"
$ws.Range("E45").Value = "  +3.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.10"
$ws.Range("E46").Value = "  +2.24%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5082"
$ws.Range("E47").Value = "  -0.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4046"
$ws.Range("E48").Value = "  +1.61%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.952"
$ws.Range("E49").Value = "  +1.36%  "

$ws.Range("B50").Value = "XinFinNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07200"
$ws.Range("E50").Value = "  -1.66%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05783"
$ws.Range("E51").Value = "  +0.42%  "

